$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 / Row 3 content swap (per diff) ---
# Row 2 becomes the "merico" row with the newer record time and a 5th email (L2)
$ws.Range("F2").Value2 = "2023-05-27 18:46:39 [UTC+08:00]"
$ws.Range("G2").Value2 = "merico"
$ws.Range("H2").Value2 = "test0@merico.dev"
$ws.Range("I2").Value2 = "test1@merico.dev"
$ws.Range("J2").Value2 = "test2@merico.dev"
$ws.Range("K2").Value2 = "test3@devchat.ai"
$ws.Range("L2").Value2 = "test3@merico.dev"

# Row 3 becomes the "covespace" row, keeps the original record time, loses the 5th email (L3)
$ws.Range("G3").Value2 = "covespace"
$ws.Range("H3").Value2 = "test0@devchat.ai"
$ws.Range("I3").Value2 = "test1@devchat.ai"
$ws.Range("J3").Value2 = "test2@devchat.ai"
$ws.Range("K3").Value2 = "test3@devchat.ai"
$ws.Range("L3").Clear()

# --- Hyperlinks: re-point them so the same targets are reachable from the new cells ---
$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("H2"), "mailto:test0@merico.dev")
$ws.Hyperlinks.Add($ws.Range("I2"), "mailto:test1@merico.dev")
$ws.Hyperlinks.Add($ws.Range("J2"), "mailto:test2@merico.dev")
$ws.Hyperlinks.Add($ws.Range("K2"), "mailto:test3@devchat.ai")
$ws.Hyperlinks.Add($ws.Range("L2"), "mailto:test3@merico.dev")
$ws.Hyperlinks.Add($ws.Range("H3"), "mailto:test0@devchat.ai")
$ws.Hyperlinks.Add($ws.Range("I3"), "mailto:test1@devchat.ai")
$ws.Hyperlinks.Add($ws.Range("J3"), "mailto:test2@devchat.ai")
$ws.Hyperlinks.Add($ws.Range("K3"), "mailto:test3@devchat.ai")

# Re-apply the Hyperlink cell style so the style index matches the canonical one
# (Hyperlinks.Add() otherwise mints a near-duplicate style entry)
$ws.Range("H2").Style = "Hyperlink"
$ws.Range("I2").Style = "Hyperlink"
$ws.Range("J2").Style = "Hyperlink"
$ws.Range("K2").Style = "Hyperlink"
$ws.Range("L2").Style = "Hyperlink"
$ws.Range("H3").Style = "Hyperlink"
$ws.Range("I3").Style = "Hyperlink"
$ws.Range("J3").Style = "Hyperlink"
$ws.Range("K3").Style = "Hyperlink"

# --- Selection change ---
$ws.Range("F9").Select()
